$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column A (Livre I) and column D (FR) for rows 3-11
for ($r = 3; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "Livre I"
    $ws.Cells.Item($r, 4).Value = "FR"
}

# Update the long-form text entries in column C (rows 7-11)
$ws.Range("C7").Value = 'De Diognète, d''éviter de perdre son temps ; de ne jamais croire aux incantations et conjurations de démons des sorciers et autres charlatants, ni à tant d’autres inventions aussi fausses. De ne pas m''attarder sur des batailles futiles ou sur d''autres puérilités ; de savoir supporter la franchise de ceux qui me parlent ; d’avoir contracté le goût de la philosophie ; d’avoir suivi les leçons de Bacchius, puis ensuite celles de Tandasis et de Marcien; d’avoir écrit des essais dès mon enfance, et de m’être fait au mode de vie grec, leurs lit de camps et leurs cuir.'
$ws.Range("C8").Value = 'De Rusticus, de m''avoir fait comprendre que devais entraîner et̀ redresser mon caractère ; de ne point me laisser aller aux engouements de la réthorie ; de ne point écrire sur les sciences spéculatives, de ne pas donner de petits sermons vaniteux, ou de ne point chercher à frapper les imaginations en m’affichant pour un homme plein d’activité ou de bienfaisance ; de rester loin de la réthorie, de la poésie et de toute affectation dans le style ; de ne pas parader en costume de cérémonie à la maison, et de ne faire des choses du même genre ; d’écrire sans aucune prétention ma correspondance, dans le genre de la lettre qu’il écrivit lui- même de Sinuesse à ma mère. Il m’a montré aussi à être toujours prêt à rappeler ou à accueillir ceux qui m’avaient chagriné ou négligé, dès le moment qu’ils étaient eux- mêmes disposés à revenir ; de toujours apporter une grande attention à mes lectures, et à ne pas me contenter de comprendre à demi ce que je lisais ; d̀e ne pas acquiescer trop vite aux propositions qui m’étaient faites. Enfin, je lui dois d’avoir connu les Entretiens d’Épictète, qu’il me prêta de sa propre bibliothèque. '
$ws.Range("C9").Value = 'D’Apollonius, j’ai appris à avoir l’esprit libre et à être ferme sans hésitation ; à ne rien chercher sauf la raison, sans en dévier un seul instant ; à conserver toujours une parfaite égalité d’âme contre les douleurs les plus vives, la perte d’un enfant ou les longues maladies. À voir qu’une même personne peut être tout plein de force de caractère et de peut être relâchée ; et qu’on peut n’être point impatient en enseignant ; il m’a montré être un homme qui regarde la moindre de ses qualités de savoir transmettre la science à autrui, comme des qualités basiques. C’est lui encore qui m’a appris l’art de recevoir de mes amis des services, sans en être diminué, et sans y paraître insensible quand je ne croyais pas devoir les accepter.'
$ws.Range("C10").Value = 'De Sextus, j’ai appris ce qu’est la bienveillance, une famille paternellement gouvernée et le vrai sens du précepte Vivre selon la Nature ; la gravité sans prétention ; de montrer de la sympathie pour nos amis ; la patience à supporter les amateurs et leurs propos irréfléchis ; la faculté de s’entendre si bien avec tout le monde que sa simple présence semblait plus agréable que ne peut l’être aucune flatterie, et que ceux qui l’entretenaient n’avaient jamais plus de respect pour lui que dans ces rencontres ; l’habileté à saisir, à trouver, et à classer les préceptes nécessaires à la pratique de la vie ; le soin de ne jamais montrer d’emportement ni aucune autre passion excessive ; le talent d’être à la fois le plus impassible et le plus affectueux des hommes ; le plaisir à dire du bien des gens mais sans l''ébruiter ; enfin d''être instruit sans le montrer.'
$ws.Range("C11").Value = 'Par l’exemple d’Alexandre le grammairien, j’ai appris à ne jamais corriger les gens, surtout par une brusquerie blessante pour un barbarisme ou une mauvaise prononciation qu’ils auraient commis ; mais à m’arranger adroitement dans la conversation pour que le mot qui aurait dû être choisi d’abord reparût, en réponse ou en confirmation, en donnant mon avis sur la chose même sans m’arrêter du tout à l’expression malheureuse, ou en prenant soigneusement tel autre détour pour dissimuler l’allusion. '

# Update selection to match the authored state
$ws.Range("E4").Select()
